$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab (Report -> Sheet1)
$ws.Name = "Sheet1"

# Update membership-count column (B) for each society row.
# Row 2 = FLASCO, Row 3 = GASCO, Row 4 = IOS, Row 5 = IOWA, Row 6 = MOASC
$ws.Range("B2").Value = 120
$ws.Range("B3").Value = 190
$ws.Range("B4").Value = 150
$ws.Range("B5").Value = 240
$ws.Range("B6").Value = 900
